# egg grading considering number of eggs
# Extends the weekly tracking table with 4 more week blocks (Week 19-22),
# recomputes the cumulative "Feed" totals in row 3 for the existing weeks,
# and fills in the data for the newly added weeks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Column widths for the 4 new week blocks (columns 82-97, CD:CS)
#    Pattern mirrors the existing blocks: first column of each block is
#    slightly wider (Weight col), the following three are narrower.
# ---------------------------------------------------------------------
$newBlockStarts = @(82, 86, 90, 94)
foreach ($start in $newBlockStarts) {
    $ws.Columns.Item($start).ColumnWidth = 6.42
    $ws.Columns.Item($start + 1).ColumnWidth = 6.25
    $ws.Columns.Item($start + 2).ColumnWidth = 6.25
    $ws.Columns.Item($start + 3).ColumnWidth = 6.25
}

# ---------------------------------------------------------------------
# 2) Row 1 - merged "Week N" header blocks, copied (format+value) from
#    the last existing block (BZ1:CC1 = "Week 18") then relabeled.
# ---------------------------------------------------------------------
$weekBlocks = @(
    @{ Dest = "CD1"; Range = "CD1:CG1"; Label = "Week 19" },
    @{ Dest = "CH1"; Range = "CH1:CK1"; Label = "Week 20" },
    @{ Dest = "CL1"; Range = "CL1:CO1"; Label = "Week 21" },
    @{ Dest = "CP1"; Range = "CP1:CS1"; Label = "Week 22" }
)
foreach ($blk in $weekBlocks) {
    $ws.Range("BZ1:CC1").Copy($ws.Range($blk.Dest))
    $ws.Range($blk.Dest).Value = $blk.Label
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Row 2 - sub headers (Weight / Egg / Egg Weight / Feed), copied
#    verbatim (format+value) from the last existing block.
# ---------------------------------------------------------------------
$ws.Range("BZ2:CC2").Copy($ws.Range("CD2"))
$ws.Range("BZ2:CC2").Copy($ws.Range("CH2"))
$ws.Range("BZ2:CC2").Copy($ws.Range("CL2"))
$ws.Range("BZ2:CC2").Copy($ws.Range("CP2"))
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4) Row 3 - data row. First update the recalculated cumulative "Feed"
#    values for the pre-existing weeks (now factoring in egg counts),
#    then copy formatting into the new blocks and set their values.
# ---------------------------------------------------------------------
$existingUpdates = @{
    "Q3" = 126; "U3" = 182; "Y3" = 231; "AC3" = 266; "AG3" = 301;
    "AK3" = 329; "AO3" = 357; "AS3" = 385; "AW3" = 413; "BA3" = 434;
    "BE3" = 455; "BI3" = 476; "BM3" = 497; "BQ3" = 518; "BU3" = 539;
    "BY3" = 567; "CA3" = 1; "CB3" = 50; "CC3" = 609
}
foreach ($ref in $existingUpdates.Keys) {
    $ws.Range($ref).Value = $existingUpdates[$ref]
}

$ws.Range("BZ3:CC3").Copy($ws.Range("CD3"))
$ws.Range("BZ3:CC3").Copy($ws.Range("CH3"))
$ws.Range("BZ3:CC3").Copy($ws.Range("CL3"))
$ws.Range("BZ3:CC3").Copy($ws.Range("CP3"))
$excel.CutCopyMode = $false

$newRow3 = @{
    "CD3" = 1560; "CE3" = 2;  "CF3" = 112; "CG3" = 621;
    "CH3" = 1610; "CI3" = 4;  "CJ3" = 260; "CK3" = 645;
    "CL3" = 1665; "CM3" = 2;  "CN3" = 152; "CO3" = 670;
    "CP3" = 1710; "CQ3" = 5;  "CR3" = 250; "CS3" = 690
}
foreach ($ref in $newRow3.Keys) {
    $ws.Range($ref).Value = $newRow3[$ref]
}

# ---------------------------------------------------------------------
# 5) Row 8 - small marker cells, one per new week block.
# ---------------------------------------------------------------------
$ws.Range("BZ8").Copy($ws.Range("CD8"))
$ws.Range("BZ8").Copy($ws.Range("CH8"))
$ws.Range("BZ8").Copy($ws.Range("CL8"))
$ws.Range("BZ8").Copy($ws.Range("CP8"))
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 6) Row 10 - bottom border row, extended across all new columns; also
#    normalise F10 (picks up the same style as its neighbours).
# ---------------------------------------------------------------------
$ws.Range("A10").Copy($ws.Range("F10"))
$ws.Range("CC10").Copy($ws.Range("CD10:CS10"))
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 7) Selection / active cell, matching the saved view state.
# ---------------------------------------------------------------------
$ws.Range("CR3").Select()
